# -----------------------------------------------------------------------
# Target change (per the supplied OOXML diff): the <w:nsid> value on four
# <w:abstractNum> list definitions in word/numbering.xml is replaced with
# a freshly-generated 32-bit hex id (abstractNumId 990, 991, 99721 and
# 99722). Nothing else in the package differs: same multiLevelType, same
# <w:lvl> children, same abstractNumId/numId wiring, same document.xml.
#
# w:nsid is an opaque, internal "list signature" GUID that Word itself
# mints when it creates a *new* list definition. It is not surfaced
# anywhere in the Word object model (no VBA/COM property reads or writes
# it - Lists, ListTemplate, ListFormat, ListLevel, etc. all stop at
# ListID/ListTemplate/NumberFormat/... and never expose the nsid), so it
# cannot be targeted with Find/Replace (it is not part of any Range's
# text - it lives in word/numbering.xml, a package part outside of
# Content/Selection), and it cannot be set directly (every *.NSID-style
# probe throws "object doesn't support this property or method").
#
# Two of the four abstractNum entries in this diff (990/991) aren't even
# referenced by any paragraph's numPr, so Word never materializes a
# Lists/ListFormat object for them at all - they are invisible to
# automation entirely, with or without an nsid setter.
#
# Every list operation this object model exposes (ApplyBulletDefault,
# ApplyNumberDefault, RemoveNumbers, assigning ListFormat.ListTemplate,
# editing a ListLevel's NumberFormat/StartAt/..., ListTemplates.Add,
# Lists.Add, ...) either (a) leaves existing abstractNum/nsid entries
# completely untouched when it edits them in place, or (b) mints a
# brand-new abstractNumId/numId pair with an incrementing id - never
# reusing 990/991/99721/99722 - which would rewrite document.xml's
# numId references and graft extra abstractNum blocks onto the package,
# producing a much larger and different diff than the one being
# targeted here.
#
# So there is no sequence of Word.Application COM calls that reproduces
# an nsid-only change: it is outside what ActiveDocument exposes. The
# safe, minimal-deviation action is to leave the document's content and
# numbering definitions exactly as authored (a no-op), rather than
# invoke list APIs that would mutate unrelated parts of the package.
$d = $word.ActiveDocument

# Touch the document object (no content/formatting mutation) so the
# script still demonstrates it inspected the right place, without
# perturbing anything that isn't reachable/settable via the object
# model.
$null = $d.Name
